# fix import with password&otority and create kurikulum with dropdown
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (new shared strings "Example_password" and
# "Otoritas (Dosen/Penjamin Mutu)") extending the used range to A1:D4.
$ws.Range("C1").Value = "Example_password"
$ws.Range("D1").Value = "Otoritas (Dosen/Penjamin Mutu)"

# Update the active selection to match the new workbook state.
$ws.Range("F3").Select() | Out-Null
